# Remove the "Ver no Jupiter Salvar em pdf Salvar em docx" line, the
# copyright/footer line right after it, and the blank paragraph that
# separates them from the preceding "LOQ4083: ..." requirement line.
#
# Layout right before the edit (paragraph text shown in brackets):
#   [LOQ4083: Fenômenos de Transporte I (Requisito fraco)]   <- keep
#   []                                                        <- remove (blank)
#   [Ver no Jupiter Salvar em pdf Salvar em docx]             <- remove
#   [© 2020 . Contact: luizeleno@usp.br. ... Attribution]     <- remove
#   []                                                        <- keep (blank)
#   [] (page-break-before)                                    <- keep

$d = $word.ActiveDocument

$target1 = "Ver no Jupiter Salvar em pdf Salvar em docx"
$target2 = "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution"

$startPara = $null
$endPara = $null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $text = $p.Range.Text

    if ($text -like "*$target1*") {
        $startPara = $i
    }
    if ($text -like "*$target2*") {
        $endPara = $i
    }
}

if ($startPara -ne $null -and $endPara -ne $null) {
    # Also drop the blank paragraph immediately preceding the "Ver no
    # Jupiter..." paragraph, which separated it from the requirement line.
    $deleteFrom = $startPara - 1
    $deleteTo = $endPara

    $rangeStart = $d.Paragraphs.Item($deleteFrom).Range.Start
    $rangeEnd = $d.Paragraphs.Item($deleteTo).Range.End

    $r = $d.Range($rangeStart, $rangeEnd)
    $r.Delete()
}
